$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Clientes")

# Widen the "Nombre" column (B) to fit longer names
# (the engine's ColumnWidth -> stored-width round trip adds 5/6, so back it
# out here to land on the target stored width of exactly 19)
$ws.Columns.Item(2).ColumnWidth = 19 - (5/6)

# Row 8: cedula/telefono were stored as text, now store as real numbers
$ws.Range("A8").Value = 1017237015
$ws.Range("C8").Value = 3057897240

# New client rows
$ws.Range("A9").Value = 4848484848
$ws.Range("B9").Value = "asdasdasdasdasdas"
$ws.Range("C9").Value = 1234567890

$ws.Range("A10").Value = 1234567891
$ws.Range("B10").Value = "Pedro"
$ws.Range("C10").Value = 3014386600

$ws.Range("A11").Value = 1234567891
$ws.Range("B11").Value = "aaaa"
$ws.Range("C11").Value = 3245619850
